$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last status check" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 14.02.2022 01:15"

# Row 5 (Makro) gets new price readings:
#   B5 (Cena)      -> new current price
#   C5 (Old Cena)  -> previous current price (old B5 value)
#   D5 (Delta)     -> text "+0.6"
#   E5 (Old Datum) -> text timestamp string
$ws.Range("B5").Value = 36.1
$ws.Range("C5").Value = 35.5

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "+0.6"
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2022-02-14 01:15:07"
$ws.Range("E5").Style = "Normal"
